# The deck ships two theme parts:
#   theme1.xml - "Integral"     (bound to the slide master -> drives every slide)
#   theme2.xml - "Office Theme" (bound only to the notes master, otherwise unused)
#
# The authored edit swaps the two themes' content, so the slide master (and
# therefore the whole deck) now renders with the plain "Office Theme" palette
# instead of "Integral". The font scheme and format scheme (fills/lines/
# effects) are already byte-identical between the two theme parts, so the
# only substantive difference to reproduce is the 12-colour theme colour
# scheme that decorates the slide master's theme.
#
# Apply the target ("Office Theme") colours onto the presentation's theme
# colour scheme via PowerPoint's ThemeColorScheme object, which is the
# documented, supported way to edit a theme's colours through the object
# model (Design tab > Colors in the UI).

function HexToComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Office Theme's 12 scheme colours, in ThemeColorScheme.Colors(index) order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToComRgb $officeThemeColors[$i - 1]
}
